# MASTER: Fix Fisher LDA.
# Adds a "Class" (music genre) column and an "AUC" metric column to the
# normalization-function comparison sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "Normalization function"
$ws.Range("B1").Value = "Class"
$ws.Range("C1").Value = "MSE"
$ws.Range("D1").Value = "Accuracy"
$ws.Range("E1").Value = "Specificity"
$ws.Range("F1").Value = "Sensitivity"
$ws.Range("G1").Value = "F-measure"
$ws.Range("H1").Value = "AUC"

# --- Column B (Class / genre) values for every data row -------------------
$genres = @("blues","classical","country","disco","hiphop","jazz","metal","pop","reggae","rock")
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 2).Value = $genres[($r - 2) % 10]
}

# --- Column A (Normalization function) values ------------------------------
for ($r = 2; $r -le 11; $r++)  { $ws.Cells.Item($r, 1).Value = "zscore" }
for ($r = 12; $r -le 21; $r++) { $ws.Cells.Item($r, 1).Value = "norm" }
for ($r = 22; $r -le 31; $r++) { $ws.Cells.Item($r, 1).Value = "range" }

# --- New column H (AUC) values ---------------------------------------------
$auc = @{
    2  = 0.48472222222222222
    3  = 0.9376535626535627
    4  = 0.44791666666666669
    5  = 0.47631578947368419
    6  = 0.49051633298208641
    7  = 0.48472222222222222
    8  = 0.47308488612836441
    9  = 0.44525547445255476
    10 = 0.44718309859154931
    11 = 0.44932432432432434
    12 = 0.50331785003317853
    13 = 0.90051020408163263
    14 = 0.40277777777777779
    15 = 0.43434134217067111
    16 = 0.46727423363711684
    17 = 0.57932330827067668
    18 = 0.55190417690417692
    19 = 0.40344827586206899
    20 = 0.400709219858156
    21 = 0.44261083743842361
    22 = 0.39930555555555558
    23 = 0.89834515366430245
    24 = 0.3971631205673759
    25 = 0.40202702702702703
    26 = 0.397887323943662
    27 = 0.59735872235872245
    28 = 0.50939597315436236
    29 = 0.4263157894736842
    30 = 0.46809440559440568
    31 = 0.49074074074074081
}

for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 8).Value = $auc[$r]
}

# --- Column widths: B widened, new H matches the other metric columns -----
# (ColumnWidth is quantized to Excel's 1/6-character grid on write, so these
# are the closest settable values to the target 8.28515625 / 12.7109375.)
$ws.Columns.Item(2).ColumnWidth = 7.5
$ws.Columns.Item(8).ColumnWidth = 11.833333333333334
